$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A ("terminus") - shifts B..E left to A..D
$ws.Columns.Item(1).Delete()

# Now columns are: A=name, B=sequence, C=Fixed_Fprimer, D=Fixed_Rprimer
# Delete columns C:D (Fixed_Fprimer, Fixed_Rprimer)
$ws.Range("C:D").Delete()

# Delete row 3's content (A3 = MBP, B3 = old long seq) but keep formatting on B3
$ws.Range("A3:B3").ClearContents()

# Update headers
$ws.Range("A1").Value = "tagName"
$ws.Range("B1").Value = "sequence"

# Update row2 values
$ws.Range("A2").Value = "MBP"
$ws.Range("B2").Value = "AAGATTGAAGAAGGTAAGTTGGTTATCTGGATTAACGGTGACAAGGGTTACAACGGTTTGGCTGAAGTTGGTAAGAAATTTGAAAAAGATACCGGTATCAAGGTCACTGTTGAACACCCAGACAAGTTGGAAGAAAAGTTTCCACAAGTTGCTGCCACTGGTGATGGTCCAGACATTATCTTCTGGGCTCATGACAGATTCGGTGGTTACGCCCAATCCGGTTTGTTAGCCGAGATCACCCCAGATAAGGCTTTTCAAGATAAGTTGTATCCATTCACTTGGGATGCCGTCAGATACAACGGTAAGTTAATCGCCTACCCAATTGCTGTTGAAGCTTTGTCTTTGATCTACAATAAGGACTTGTTACCTAACCCACCAAAGACCTGGGAAGAAATCCCAGCTTTAGATAAGGAGTTAAAAGCTAAGGGTAAGTCCGCTTTGATGTTTAACTTGCAAGAACCATACTTCACTTGGCCATTGATCGCTGCTGATGGTGGTTACGCTTTTAAGTATGAAAACGGTAAATACGACATTAAGGATGTCGGTGTCGACAATGCTGGTGCTAAGGCCGGTTTAACTTTCTTAGTCGATTTGATTAAGAATAAACATATGAATGCTGACACTGATTACTCTATTGCTGAAGCTGCTTTCAACAAGGGTGAAACCGCTATGACTATTAACGGTCCATGGGCCTGGTCTAACATTGATACCTCTAAAGTCAACTACGGTGTCACCGTCTTGCCAACTTTTAAGGGTCAACCATCTAAGCCATTCGTCGGTGTCTTGTCTGCCGGTATTAACGCTGCCTCTCCAAATAAGGAATTGGCCAAGGAATTCTTAGAAAACTACTTGTTAACCGATGAAGGTTTAGAGGCCGTTAACAAGGATAAGCCATTAGGTGCTGTTGCTTTGAAGTCTTACGAAGAAGAGTTGGCTAAGGATCCAAGAATTGCTGCTACTATGGAAAACGCTCAAAAGGGTGAAATTATGCCAAACATCCCACAAATGTCTGCTTTCTGGTACGCTGTTCGTACCGCCGTCATTAATGCCGCTTCTGGTCGTCAAACTGTTGATGAAGCCTTGAAGGACGCTCAAACCAGAATTACTAAG"

# Adjust row 3 height
$ws.Rows.Item(3).RowHeight = 16.5

# Selection
$ws.Range("B3").Select()

Write-Output $ws.UsedRange.Address()
